$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode column (D) to "N" for all rows except row 12 (Notifications012),
# which remains "Y" so that only that one notification scenario runs.
for ($r = 2; $r -le 26; $r++) {
    if ($r -ne 12) {
        $ws.Cells.Item($r, 4).Value = "N"
    }
}

# Update the active selection on the sheet.
$ws.Range("C18").Select()
